$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.489.99'
$ws.Range("E2").Value = '  +3.62%  '
$ws.Range("D3").Value = '3.500.16'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'590.64"
$ws.Range("E5").Value = '  +3.23%  '
$ws.Range("D6").Value = "'169.80"
$ws.Range("E6").Value = '  +5.32%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.497.39'
$ws.Range("E8").Value = '  +2.11%  '
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = '  +6.25%  '
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("E11").Value = '  +4.86%  '
$ws.Range("E12").Value = '  +3.34%  '
$ws.Range("D13").Value = '4.104.70'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = "'28.31"
$ws.Range("E15").Value = '  +4.79%  '
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("D17").Value = '66.500.88'
$ws.Range("E17").Value = '  +3.55%  '
$ws.Range("D18").Value = '3.502.81'
$ws.Range("E18").Value = '  +2.84%  '
$ws.Range("E19").Value = '  +3.89%  '
$ws.Range("D20").Value = "'14.10"
$ws.Range("E20").Value = '  +3.82%  '
$ws.Range("D21").Value = "'390.89"
$ws.Range("E21").Value = '  +3.74%  '
$ws.Range("D22").Value = "'7.99"
$ws.Range("E22").Value = '  +2.34%  '
$ws.Range("D23").Value = "'72.97"
$ws.Range("E23").Value = '  +2.35%  '
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").Value = "'0.536"
$ws.Range("E25").Value = '  +3.45%  '
$ws.Range("E26").Value = '  +5.02%  '
$ws.Range("D27").Value = "'10.46"
$ws.Range("E27").Value = '  +9.51%  '
$ws.Range("E28").Value = '  +2.45%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").Value = "'6.35"
$ws.Range("E30").Value = '  +5.62%  '
$ws.Range("D32").Value = "'2.07"
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("D33").Value = "'23.62"
$ws.Range("E33").Value = '  +3.08%  '
$ws.Range("D34").Value = "'7.42"
$ws.Range("E34").Value = '  +4.28%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +7.62%  '
$ws.Range("D37").Value = "'162.37"
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("D38").Value = "'0.883"
$ws.Range("E38").Value = '  +3.06%  '
$ws.Range("E39").Value = '  +5.06%  '
$ws.Range("D40").Value = "'6.86"
$ws.Range("E40").Value = '  +6.28%  '
$ws.Range("D41").Value = "'27.84"
$ws.Range("E41").Value = '  +6.95%  '
$ws.Range("E42").Value = '  +6.27%  '
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("D44").Value = "'26.52"
$ws.Range("E44").Value = '  +2.84%  '
$ws.Range("D45").Value = '2.794.15'
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").Value = "'0.0312"
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("D48").Value = "'2.50"
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = "'351.92"
$ws.Range("E49").Value = '  +5.06%  '
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("D51").Value = "'33.70"
$ws.Range("E51").Value = '  +12.17%  '
